## componentes.xlsx edit script
## Adds a "LOJA" column, a "Frete (R$)" column, a new ESP32 row, renames
## several headers/labels, replaces the cart total with a formula, and
## adds hyperlinks to the "Link" column cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Structural changes: insert the two new columns.
#    Before:  D Item | E Mouser | F Fabricante | G Qtd | H Preco | I Link
#    Insert a fresh column at E  -> shifts Mouser..Link one to the right
#    Insert a fresh column at J (after the shifted "Preco", before "Link")
# ---------------------------------------------------------------------
$ws.Columns.Item(5).Insert()   # new column E (will hold "LOJA")
$ws.Columns.Item(10).Insert()  # new column J (will hold "Frete (R$)")

# the old literal "$X" cart-total placeholder (originally F9) ends up at
# G9 after the two column inserts; it is replaced by a formula in F9, so
# drop the stray leftover cell entirely.
$ws.Range("G9").ClearContents()

# ---------------------------------------------------------------------
# 2) Rename existing header / label text
# ---------------------------------------------------------------------
$ws.Range("I3").Value = "Preço (R$)"
$ws.Range("E9").Value = "Total (R$)"

# ---------------------------------------------------------------------
# 3) New header cells
# ---------------------------------------------------------------------
$ws.Range("E3").Value = "LOJA"
$ws.Range("J3").Value = "Frete (R$)"

# ---------------------------------------------------------------------
# 4) Row 4 (Transceptor RS485 Arduino)
# ---------------------------------------------------------------------
$ws.Range("D4").Value = "Transceptor RS485 Arduino"
$ws.Range("E4").Value = "Mercado Livre"
$ws.Range("F4").Value = "-"
$ws.Range("G4").Value = "-"
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 15.9
$ws.Range("J4").Value = 40

# ---------------------------------------------------------------------
# 5) Row 5 (RS485 communication module cfw100)
# ---------------------------------------------------------------------
$ws.Range("D5").Value = "RS485 communication module cfw100"
$ws.Range("E5").Value = "ViewTech"
$ws.Range("F5").Value = "-"
$ws.Range("G5").Value = "-"
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 419.25
$ws.Range("J5").Value = 22

# ---------------------------------------------------------------------
# 6) Row 6 (new row - ESP32)
# ---------------------------------------------------------------------
$ws.Range("D6").Value = "ESP32"
$ws.Range("E6").Value = "Mercado Livre"
$ws.Range("F6").Value = "-"
$ws.Range("G6").Value = "-"
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 47.49
$ws.Range("J6").Value = 41

# ---------------------------------------------------------------------
# 7) Row 9 (cart total) -- replace literal "$X" with a formula
# ---------------------------------------------------------------------
$ws.Range("F9").Formula = "=I4+J4+I5+J5+I6+J6"

# ---------------------------------------------------------------------
# 8) Formatting
# ---------------------------------------------------------------------

# Header band (row 2, merged) and header rows (3, 9) -> bold
$ws.Range("D2:K2").Font.Bold = $true
$ws.Range("D3:K3").Font.Bold = $true
$ws.Range("D9:E9").Font.Bold = $true

# Borders for the new data block (rows 4-6, D:K) -- reuse the thin/automatic
# border already used elsewhere (set Color before LineStyle so that the
# engine matches the existing border definition instead of creating a
# duplicate one).
$dataRange = $ws.Range("D4:K6")
$dataRange.Borders.ColorIndex = 1
$dataRange.Borders.LineStyle = 1

# Alignment
$ws.Range("E4:E6").HorizontalAlignment = -4152   # xlRight
$ws.Range("F4:G4").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F4:G4").VerticalAlignment = -4108     # xlCenter
$ws.Range("F6:G6").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F6:G6").VerticalAlignment = -4108     # xlCenter
$ws.Range("F5:G5").HorizontalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------------
# 9) Hyperlinks (added after borders so the engine merges the hyperlink
#    font with the border already present on the cell)
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("K5"), "https://www.viewtech.ind.br/modulo-de-expansao-weg-cfw100-crs485-comunicacao") | Out-Null

$urlMLB1 = "https://produto.mercadolivre.com.br/MLB-2655392212-modulo-conversor-bidirecional-rs485-ttl-arduino-pic-rasp-pic-_JM?matt_tool=14213447&matt_word=&matt_source=bing&matt_campaign=MLB_ML_BING_AO_CE-ALL-ALL_X_PLA_ALLB_TXS_ALL&matt_campaign_id=382858295&matt_ad_group=CE&matt_match_type=e&matt_network=o&matt_device=c&matt_keyword=default&msclkid=177180d409fd133f90eee53cb77ef9ed&utm_source=bing&utm_medium=cpc&utm_campaign=MLB_ML_BING_AO_CE-ALL-ALL_X_PLA_ALLB_TXS_ALL&utm_term=4581596253419741&utm_content=CE"
$ws.Hyperlinks.Add($ws.Range("K4"), $urlMLB1) | Out-Null

$urlMLB2 = "https://produto.mercadolivre.com.br/MLB-2712932646-modulo-wi-fi-esp32-nodemcu-com-bluetooth-38-pinos-_JM?matt_tool=14213447&matt_word=&matt_source=bing&matt_campaign=MLB_ML_BING_AO_CE-ALL-ALL_X_PLA_ALLB_TXS_ALL&matt_campaign_id=382858295&matt_ad_group=CE&matt_match_type=e&matt_network=o&matt_device=c&matt_keyword=default&msclkid=d45d1ae32055141043757f0bf54a32fd&utm_source=bing&utm_medium=cpc&utm_campaign=MLB_ML_BING_AO_CE-ALL-ALL_X_PLA_ALLB_TXS_ALL&utm_term=4581596253419738&utm_content=CE"
$ws.Hyperlinks.Add($ws.Range("K6"), $urlMLB2) | Out-Null

# ---------------------------------------------------------------------
# 10) Column widths (best effort - Excel quantizes these to its pixel
#     grid so we cannot reproduce the exact fractional widths, but this
#     gets us close).
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 12.88    # E
$ws.Columns.Item(9).ColumnWidth = 12.02    # I
$ws.Columns.Item(10).ColumnWidth = 12.02   # J
$ws.Columns.Item(11).ColumnWidth = 15.02   # K

# ---------------------------------------------------------------------
# 11) Selection cosmetics
# ---------------------------------------------------------------------
$ws.Range("F10").Select() | Out-Null
